$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.961.68"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").Value = "2.466.05"
$ws.Range("E3").Value = "  -2.35%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'518.16"
$ws.Range("E5").Value = "  -3.56%  "
$ws.Range("D6").Value = "'130.92"
$ws.Range("E6").Value = "  -4.50%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -2.26%  "
$ws.Range("D9").Value = "'0.0990"
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D13").Value = "2.901.38"
$ws.Range("E13").Value = "  -1.88%  "
$ws.Range("D14").Value = "57.860.56"
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("D15").Value = "'22.26"
$ws.Range("E15").Value = "  -3.35%  "
$ws.Range("E16").Value = "  -2.54%  "
$ws.Range("D17").Value = "2.462.84"
$ws.Range("E17").Value = "  -2.60%  "
$ws.Range("D18").Value = "'10.76"
$ws.Range("E18").Value = "  -3.67%  "
$ws.Range("E19").Value = "  -2.63%  "
$ws.Range("D20").Value = "'319.83"
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("E22").Value = "  -3.75%  "
$ws.Range("D23").Value = "'64.02"
$ws.Range("E23").Value = "  -2.91%  "
$ws.Range("D24").Value = "'0.410"
$ws.Range("E24").Value = "  -3.17%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "'0.160"
$ws.Range("E26").Value = "  -3.13%  "
$ws.Range("D27").Value = "'7.31"
$ws.Range("E27").Value = "  -3.13%  "
$ws.Range("D28").Value = "0.0₃0751"
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'1.69"
$ws.Range("E29").Value = "  -4.66%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'165.80"
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("E31").Value = "  -5.93%  "
$ws.Range("E32").Value = "  -2.63%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").Value = "'18.04"
$ws.Range("E35").Value = "  -2.12%  "
$ws.Range("E36").Value = "  -10.84%  "
$ws.Range("E37").Value = "  -3.40%  "
$ws.Range("E38").Value = "  -4.58%  "
$ws.Range("D39").Value = "'0.788"
$ws.Range("E39").Value = "  -3.28%  "
$ws.Range("D40").Value = "'3.45"
$ws.Range("E40").Value = "  -4.65%  "
$ws.Range("D41").Value = "'271.89"
$ws.Range("E41").Value = "  -4.46%  "
$ws.Range("D42").Value = "'4.99"
$ws.Range("E42").Value = "  -2.94%  "
$ws.Range("D43").Value = "'0.590"
$ws.Range("E43").Value = "  -2.79%  "
$ws.Range("D44").Value = "'126.09"
$ws.Range("E44").Value = "  -4.86%  "
$ws.Range("E45").Value = "  -2.26%  "
$ws.Range("E46").Value = "  -4.17%  "
$ws.Range("E47").Value = "  -3.45%  "
$ws.Range("D48").Value = "'17.02"
$ws.Range("E48").Value = "  -1.95%  "
$ws.Range("D49").Value = "1.729.37"
$ws.Range("E49").Value = "  -2.13%  "
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("D51").Value = "'4.68"
$ws.Range("E51").Value = "  -1.29%  "
